$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.3
$ws.Range("B21").Value = 9.570000000000002
$ws.Range("B23").Value = 7.43
$ws.Range("E24").Value = 16.803
$ws.Range("B25").Value = 6.396
$ws.Range("C27").Value = -13.074
$ws.Range("C31").Value = -13.163
$ws.Range("C39").Value = -12.858
$ws.Range("C48").Value = -11.432
$ws.Range("C51").Value = -11.133
$ws.Range("C52").Value = -11.273
$ws.Range("B53").Value = 5.697
$ws.Range("C55").Value = -13.41
$ws.Range("C56").Value = -13.537
$ws.Range("B57").Value = 4.864000000000001
$ws.Range("C57").Value = -13.852
$ws.Range("E57").Value = 16.334
$ws.Range("B59").Value = 4.678
$ws.Range("E61").Value = 16.706
$ws.Range("B69").Value = 5.667
$ws.Range("E70").Value = 17.687
$ws.Range("C73").Value = -12.575
$ws.Range("B79").Value = 5.855
$ws.Range("B83").Value = 5.702
$ws.Range("E86").Value = 16.554
$ws.Range("C89").Value = -11.418
$ws.Range("C90").Value = -12.91
$ws.Range("B93").Value = 5.659000000000001
$ws.Range("E98").Value = 16.239
$ws.Range("E100").Value = 16.809
$ws.Range("E102").Value = 16.536
